# Revert "feat(dialog): update CN data and dialogue Excel files"
# - Restore the hotspring_meditation row's version/text/text_EN/text_JP to
#   their pre-update ("meditation") values.
# - Remove the pond_carp row that was added by the reverted commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Note")

# Update row 13 (hotspring_meditation) back to its original values.
$ws.Range("B13").Value = "EA 23.27 fix 2"
$ws.Range("C13").Value = "在温泉中冥想，可以温暖身心。"
$ws.Range("D13").Value = "Hot spring meditation will warm you from the inside out!"
$ws.Range("E13").Value = "温泉で瞑想すれば、心も体もほっこり"

# Remove the entire pond_carp row (row 14).
$ws.Rows.Item(14).Delete()
